$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: Volume Number and Report Week dates ---
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Brooklyn North precinct weekly crime data (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = "0"
$ws.Range("E14").Value = "***.*"
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = 50
$ws.Range("I14").Value = 24
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 20
$ws.Range("L14").Value = -20
$ws.Range("M14").Value = -42.857142857142
$ws.Range("N14").Value = -85.454545454545

# Row 15
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 16
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 83
$ws.Range("J15").Value = 85
$ws.Range("K15").Value = -2.35294117647
$ws.Range("L15").Value = 7.792207792207
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -63.755458515283

# Row 16
$ws.Range("C16").Value = 49
$ws.Range("D16").Value = 60
$ws.Range("E16").Value = -18.333333333333
$ws.Range("F16").Value = 190
$ws.Range("G16").Value = 208
$ws.Range("H16").Value = -8.653846153846
$ws.Range("I16").Value = 851
$ws.Range("J16").Value = 895
$ws.Range("K16").Value = -4.916201117318
$ws.Range("L16").Value = 28.355957767722
$ws.Range("M16").Value = -27.077977720651
$ws.Range("N16").Value = -85.654079568442

# Row 17
$ws.Range("C17").Value = 83
$ws.Range("D17").Value = 91
$ws.Range("E17").Value = -8.791208791208
$ws.Range("F17").Value = 319
$ws.Range("G17").Value = 307
$ws.Range("H17").Value = 3.908794788273
$ws.Range("I17").Value = 1440
$ws.Range("J17").Value = 1396
$ws.Range("K17").Value = 3.151862464183
$ws.Range("L17").Value = 29.263913824057
$ws.Range("M17").Value = 29.963898916967
$ws.Range("N17").Value = -50.76923076923

# Row 18
$ws.Range("C18").Value = 31
$ws.Range("D18").Value = 60
$ws.Range("E18").Value = -48.333333333333
$ws.Range("F18").Value = 145
$ws.Range("G18").Value = 205
$ws.Range("H18").Value = -29.268292682926
$ws.Range("I18").Value = 752
$ws.Range("J18").Value = 907
$ws.Range("K18").Value = -17.089305402425
$ws.Range("L18").Value = 15.514592933947
$ws.Range("M18").Value = -23.732251521298
$ws.Range("N18").Value = -81.875150638708

# Row 19
$ws.Range("C19").Value = 106
$ws.Range("D19").Value = 109
$ws.Range("E19").Value = -2.752293577981
$ws.Range("F19").Value = 431
$ws.Range("G19").Value = 454
$ws.Range("H19").Value = -5.066079295154
$ws.Range("I19").Value = 2040
$ws.Range("J19").Value = 1986
$ws.Range("K19").Value = 2.719033232628
$ws.Range("L19").Value = 38.586956521739
$ws.Range("M19").Value = 48.68804664723
$ws.Range("N19").Value = -9.534368070953

# Row 20
$ws.Range("C20").Value = 45
$ws.Range("D20").Value = 42
$ws.Range("E20").Value = 7.142857142857
$ws.Range("F20").Value = 137
$ws.Range("G20").Value = 131
$ws.Range("H20").Value = 4.580152671755
$ws.Range("I20").Value = 611
$ws.Range("J20").Value = 635
$ws.Range("K20").Value = -3.779527559055
$ws.Range("L20").Value = 34.878587196468
$ws.Range("M20").Value = 27.824267782426
$ws.Range("N20").Value = -82.808103545301

# Row 21
$ws.Range("C21").Value = 321
$ws.Range("D21").Value = 367
$ws.Range("E21").Value = -12.534059945504
$ws.Range("F21").Value = 1244
$ws.Range("G21").Value = 1325
$ws.Range("H21").Value = -6.113207547169
$ws.Range("I21").Value = 5801
$ws.Range("J21").Value = 5924
$ws.Range("K21").Value = -2.076299797434
$ws.Range("L21").Value = 30.067264573991
$ws.Range("M21").Value = 10.790679908327
$ws.Range("N21").Value = -69.800614295382

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 9
$ws.Range("E22").Value = -44.444444444444
$ws.Range("F22").Value = 23
$ws.Range("G22").Value = 28
$ws.Range("H22").Value = -17.857142857142
$ws.Range("I22").Value = 114
$ws.Range("J22").Value = 140
$ws.Range("K22").Value = -18.571428571428
$ws.Range("L22").Value = 29.545454545454
$ws.Range("M22").Value = -21.917808219178
$ws.Range("N22").Value = "***.*"

# Row 23
$ws.Range("C23").Value = 29
$ws.Range("D23").Value = 37
$ws.Range("E23").Value = -21.621621621621
$ws.Range("F23").Value = 110
$ws.Range("G23").Value = 125
$ws.Range("H23").Value = -12
$ws.Range("I23").Value = 577
$ws.Range("J23").Value = 530
$ws.Range("K23").Value = 8.867924528301
$ws.Range("L23").Value = 16.096579476861
$ws.Range("M23").Value = 56.793478260869
$ws.Range("N23").Value = "***.*"

# Row 24
$ws.Range("C24").Value = 220
$ws.Range("D24").Value = 272
$ws.Range("E24").Value = -19.117647058823
$ws.Range("F24").Value = 949
$ws.Range("G24").Value = 1025
$ws.Range("H24").Value = -7.414634146341
$ws.Range("I24").Value = 4406
$ws.Range("J24").Value = 4541
$ws.Range("K24").Value = -2.972913455186
$ws.Range("L24").Value = 26.318807339449
$ws.Range("M24").Value = 29.359953024075
$ws.Range("N24").Value = "***.*"

# Row 25
$ws.Range("C25").Value = 126
$ws.Range("D25").Value = 133
$ws.Range("E25").Value = -5.263157894736
$ws.Range("F25").Value = 481
$ws.Range("G25").Value = 476
$ws.Range("H25").Value = 1.050420168067
$ws.Range("I25").Value = 2129
$ws.Range("J25").Value = 2150
$ws.Range("K25").Value = -0.976744186046
$ws.Range("L25").Value = 43.754220121539
$ws.Range("M25").Value = -22.35594456601
$ws.Range("N25").Value = "***.*"

# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = -8.333333333333
$ws.Range("I26").Value = 118
$ws.Range("J26").Value = 131
$ws.Range("K26").Value = -9.923664122137
$ws.Range("L26").Value = -11.278195488721
$ws.Range("M26").Value = "***.*"
$ws.Range("N26").Value = "***.*"

# Row 27
$ws.Range("C27").Value = 12
$ws.Range("D27").Value = 16
$ws.Range("E27").Value = -25
$ws.Range("F27").Value = 53
$ws.Range("G27").Value = 51
$ws.Range("H27").Value = 3.92156862745
$ws.Range("I27").Value = 219
$ws.Range("J27").Value = 213
$ws.Range("K27").Value = 2.81690140845
$ws.Range("L27").Value = -7.59493670886
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -20
$ws.Range("F28").Value = 17
$ws.Range("G28").Value = 23
$ws.Range("H28").Value = -26.086956521739
$ws.Range("I28").Value = 82
$ws.Range("J28").Value = 92
$ws.Range("K28").Value = -10.869565217391
$ws.Range("L28").Value = -35.433070866141
$ws.Range("M28").Value = -46.753246753246
$ws.Range("N28").Value = -88.20143884892

# Row 29
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 15
$ws.Range("G29").Value = 19
$ws.Range("H29").Value = -21.052631578947
$ws.Range("I29").Value = 70
$ws.Range("J29").Value = 81
$ws.Range("K29").Value = -13.58024691358
$ws.Range("L29").Value = -36.936936936936
$ws.Range("M29").Value = -40.17094017094
$ws.Range("N29").Value = -88.958990536277

# Row 30
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 25
$ws.Range("I30").Value = 29
$ws.Range("J30").Value = 25
$ws.Range("K30").Value = 16
$ws.Range("L30").Value = 45
$ws.Range("M30").Value = "***.*"
$ws.Range("N30").Value = "***.*"

